$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F2").Value = -4
$ws.Range("F3").Value = 1
$ws.Range("F4").Value = 2
$ws.Range("F7").Value = 9
$ws.Range("F8").Value = -8
$ws.Range("F9").Value = -7
$ws.Range("F10").Value = -1
$ws.Range("F11").Value = 1
